# Generate Report for Handoff
# - Swap the two tracked files (240b533f..., 6e6556f0...) between row 2 and row 3
#   across the Overview / zh-cn / de-de sheets.
# - The file that lands in row 3 (240b533f...) now shows as "Ready for handoff"
#   (priority "mt") with a fresh handoff timestamp; the file that lands in row 2
#   (6e6556f0...) keeps the previous "In Translation" (priority "ht") data.
# - The hyperlink targets (rId2 -> 240b533f URL, rId3 -> 6e6556f0 URL) stay put;
#   only their displayed text changes together with the cell values.
# - A couple of columns get a bit wider.

$wb = $excel.ActiveWorkbook

$urlBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob/64fb697d412e2643412a263c9a49d05ea82b1ffb/e2e/"
$url240b = $urlBase + "240b533f-824a-4f3d-9dfd-1557da7f43d1.md"
$url6e65 = $urlBase + "6e6556f0-5a2d-4cea-9731-567fb31b45d1.md"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "6e6556f0-5a2d-4cea-9731-567fb31b45d1.md"
$wsOverview.Range("A3").Value = "240b533f-824a-4f3d-9dfd-1557da7f43d1.md"

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-12 22:17:17"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $url240b, "", "", "e2e\6e6556f0-5a2d-4cea-9731-567fb31b45d1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $url6e65, "", "", "e2e\240b533f-824a-4f3d-9dfd-1557da7f43d1.md")

$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "6e6556f0-5a2d-4cea-9731-567fb31b45d1.md"
$wsZhCn.Range("G2").Value = "6e6556f0-5a2d-4cea-9731-567fb31b45d1.99bc6be93f3b393ddcc02243b2cf217c1613a329.zh-cn.xlf"

$wsZhCn.Range("A3").Value = "240b533f-824a-4f3d-9dfd-1557da7f43d1.md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("G3").Value = "240b533f-824a-4f3d-9dfd-1557da7f43d1.5bb2092c3f9eccbf59ec6a72409272fcd57d6a10.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-12 22:17:10"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $url240b, "", "", "6e6556f0-5a2d-4cea-9731-567fb31b45d1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $url6e65, "", "", "240b533f-824a-4f3d-9dfd-1557da7f43d1.md")

$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "6e6556f0-5a2d-4cea-9731-567fb31b45d1.md"
$wsDeDe.Range("G2").Value = "6e6556f0-5a2d-4cea-9731-567fb31b45d1.99bc6be93f3b393ddcc02243b2cf217c1613a329.de-de.xlf"

$wsDeDe.Range("A3").Value = "240b533f-824a-4f3d-9dfd-1557da7f43d1.md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("G3").Value = "240b533f-824a-4f3d-9dfd-1557da7f43d1.5bb2092c3f9eccbf59ec6a72409272fcd57d6a10.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-12 22:17:17"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $url240b, "", "", "6e6556f0-5a2d-4cea-9731-567fb31b45d1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $url6e65, "", "", "240b533f-824a-4f3d-9dfd-1557da7f43d1.md")

$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
